# Adds Denmark, Sweden and Norway market test-data worksheets to the
# FC_Gallery_MPM_AttachedFunctionality_FC600_Series_Panels workbook.
#
# Each new sheet is a copy of the "UK" sheet (same layout/merges/styles),
# inserted after the previous last sheet, renamed, and given its own
# "<Country> Market" / "NGC-xxxx/xxxx" values in B2/B4. Norway ends up the
# active (selected) sheet, matching the authored workbook state.

$wb = $excel.ActiveWorkbook

$uk = $wb.Worksheets.Item("UK")
$belgium = $wb.Worksheets.Item("Belgium")

# --- Denmark -----------------------------------------------------------
$uk.Copy($null, $belgium)
$denmark = $wb.Worksheets.Item($belgium.Index + 1)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2008"
$denmark.Cells.Select()

# --- Sweden --------------------------------------------------------------
$uk.Copy($null, $denmark)
$sweden = $wb.Worksheets.Item($denmark.Index + 1)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2023"
$sweden.Columns("B").ColumnWidth = 26.3
$sweden.Cells.Select()

# --- Norway --------------------------------------------------------------
$uk.Copy($null, $sweden)
$norway = $wb.Worksheets.Item($sweden.Index + 1)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1920"
$norway.Columns("B").ColumnWidth = 26.3

# Norway is the sheet left selected/active when the workbook was saved.
$norway.Activate()
$norway.Range("B2:B4").Select()
